$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 "<india>" -> "<delete>", C2 56 -> 51
$ws.Range("B2").Value = "<delete>"
$ws.Range("C2").Value = 51

# Row 3: B3 "<gape>" -> "<sea>", C3 50 -> 48
$ws.Range("B3").Value = "<sea>"
$ws.Range("C3").Value = 48

# Row 4: B4 "<sene>" -> "<she>", C4 54 -> 48
$ws.Range("B4").Value = "<she>"
$ws.Range("C4").Value = 48

# Row 5: C5 53 -> 44
$ws.Range("C5").Value = 44

# Row 6: B6 "<its>" -> "<it>"
$ws.Range("B6").Value = "<it>"

# Row 7: B7 "<whiskey>" -> "<which>"
$ws.Range("B7").Value = "<which>"

# Row 8: C8 45 -> 39
$ws.Range("C8").Value = 39
